$d = $word.ActiveDocument

function Get-ParaByText {
    param($doc, $needle)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Set-ParaText {
    param($para, $newText)
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# ---------------------------------------------------------------------------
# 1. Resolve the merge-conflict markers left in the "Results" section down to
#    the single surviving (HEAD) paragraph text.
# ---------------------------------------------------------------------------
$conflictPara = Get-ParaByText $d "<<<<<<< HEAD"
Set-ParaText $conflictPara "Using the average temperature for the upper two meters as the response variable, we initially began the variable selection process with 16 predictor variables. The variable selection process identified a reduced model with 7 variables (Figure 1). The selected variables were the average ambient air temperature for the sample date, sample date, longitude, average ambient air temperature for 30 day proceeding the sample date, elevation, latitude, length of lake shoreline, and the lake surface area."

# ---------------------------------------------------------------------------
# 2. Update the opening paragraph of "Discussion and conclusions".
# ---------------------------------------------------------------------------
$concPara = Get-ParaByText $d "Using the 2007 and 2012 NLA data"
Set-ParaText $concPara "Using the 2007 and 2012 NLA data, we have successfully built a simple yet robust model of lake surface temperature for the conterminous United States. The final model has a mean square error of 2.19 and an adjusted R^2 of 0.88. The sampling date, that day’s average ambient air temperature, data obtained from the PRISM Climate Group, and longitude are the most important variables impacting the final model’s accuracy. Given the importance of temperature to a lake ecosystem, especially cyanobacteria bloom dynamics, this model can be a valuable tool for researchers and lake resource managers. Daily predicted lake photic zone temperature for all lakes in the conterminous US can now be straightforwardly estimated based on basic ambient temperature and location information."

# ---------------------------------------------------------------------------
# 3. Update the land-use/land-cover paragraph.
# ---------------------------------------------------------------------------
$landUsePara = Get-ParaByText $d "In addition to the several average air temperature variables"
Set-ParaText $landUsePara "In addition to the several derived air temperature variables, we included land-use/land cover variables in our initial variable selection process. Specifically, we calculated the percent impervious surface for a 3km lake buffer and a measure of shoreline development (need to look at Jeff’s package to see how this is measured). This variables were included based on the hypothesis that higher amounts of development and therefore impervious surface surrounding a lake would lead to higher temperatures in lakes. Yet neither of these variables were selected in the final model. Even though the land-use variables were not selected that does not mean development and impervious are not impactful. This urban-heat effect on lakes may have been adequately captured in the average ambient air temperature. Therefore, making the land-use variables redundant. Regardless these variables did not independently contribute to the model’s accuracy."

# ---------------------------------------------------------------------------
# 4. Insert two new BodyText paragraphs right after the land-use paragraph.
# ---------------------------------------------------------------------------
$landUsePara = Get-ParaByText $d "making the land-use variables redundant"
$landUsePara.Range.InsertParagraphAfter()

$shortcomingsPara = Get-ParaByText $d "making the land-use variables redundant"
$shortcomingsPara = $shortcomingsPara.Next()
Set-ParaText $shortcomingsPara "Insert paragraph on model shortcomings"

$shortcomingsPara.Range.InsertParagraphAfter()
$sharmaPara = $shortcomingsPara.Next()
Set-ParaText $sharmaPara "Despite being one of the most common measurements collected by limnologists, lake temperature datasets that cover long periods of time are very difficult to obtain. Sharma et al (2015) have compiled summer lake temperature data for 291 lakes for the period 1985-2009. This may be the largest lake temperature database to date, however, the data are only available to members of their research group and, realistically, the number of lakes included is very small. One of the reasons we chose to model lake photic zone temperature was to develop a database of lake temperatures for the 48 conterminous United States. The model we present has proven to be accurate and will allow us to backcast lake temperatures for all the > 300,000 lakes included in NHDplus for the period of time covered by the PRISM climate predictions (1981 to present). This dataset will allow us to investigate how photic zone temperatures vary both spatially and temporally across the United States. This database is being developed and, when complete, will be made available as an open source data set."

Write-Output "done"
